$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.472738
$ws.Range("H2").Value = 1.418214
$ws.Range("I2").Value = 0.0327564895931267
$ws.Range("J2").Value = 0.03397138804734427
$ws.Range("M2").Value = 1.378421333333333
$ws.Range("N2").Value = 4.135264
$ws.Range("O2").Value = 0.01656231489052403
$ws.Range("P2").Value = 0.01794267551419991
$ws.Range("Q2").Value = 0.6516321442773335
$ws.Range("R2").Value = 5.864689298496001
$ws.Range("S2").Value = 0.0005425232953495378
$ws.Range("T2").Value = 0.0006095375925004677
$ws.Range("G3").Value = 0.472738
$ws.Range("H3").Value = 1.418214
$ws.Range("I3").Value = 0.0327564895931267
$ws.Range("J3").Value = 0.03397138804734427
$ws.Range("O3").Value = 0.2170932623988173
$ws.Range("P3").Value = 0.2351865659654651
$ws.Range("Q3").Value = 8.541375346392002
$ws.Range("R3").Value = 76.87237811752802
$ws.Range("S3").Value = 0.007111213190504782
$ws.Range("T3").Value = 0.007989614095935148
$ws.Range("G4").Value = 0.472738
$ws.Range("H4").Value = 1.418214
$ws.Range("I4").Value = 0.0327564895931267
$ws.Range("J4").Value = 0.03397138804734427
$ws.Range("M4").Value = 17.58286933333333
$ws.Range("N4").Value = 52.748608
$ws.Range("O4").Value = 0.2112656061941426
$ws.Range("P4").Value = 0.22887321273073
$ws.Range("Q4").Value = 8.312090482901333
$ws.Range("R4").Value = 74.808814346112
$ws.Range("S4").Value = 0.006920319630684036
$ws.Range("T4").Value = 0.007775140723318005
$ws.Range("G5").Value = 0.472738
$ws.Range("H5").Value = 1.418214
$ws.Range("I5").Value = 0.0327564895931267
$ws.Range("J5").Value = 0.03397138804734427
$ws.Range("M5").Value = 19.2082395
$ws.Range("N5").Value = 38.416479
$ws.Range("O5").Value = 0.2307951156866419
$ws.Range("P5").Value = 0.1666869194070983
$ws.Range("Q5").Value = 9.080464724751002
$ws.Range("R5").Value = 54.48278834850601
$ws.Range("S5").Value = 0.007560037805133958
$ws.Range("T5").Value = 0.005662586021594938
$ws.Range("G6").Value = 0.472738
$ws.Range("H6").Value = 1.418214
$ws.Range("I6").Value = 0.0327564895931267
$ws.Range("J6").Value = 0.03397138804734427
$ws.Range("M6").Value = 26.988955
$ws.Range("N6").Value = 80.966865
$ws.Range("O6").Value = 0.3242837008298742
$ws.Range("P6").Value = 0.3513106263825066
$ws.Range("Q6").Value = 12.75870460879
$ws.Range("R6").Value = 114.82834147911
$ws.Range("S6").Value = 0.01062239567145439
$ws.Range("T6").Value = 0.01193450961399571
$ws.Range("I7").Value = 0.822180234441485
$ws.Range("J7").Value = 0.8526739017519405
$ws.Range("M7").Value = 1.378421333333333
$ws.Range("N7").Value = 4.135264
$ws.Range("O7").Value = 0.01656231489052403
$ws.Range("P7").Value = 0.01794267551419991
$ws.Range("Q7").Value = 16.35581455175111
$ws.Range("R7").Value = 147.20233096576
$ws.Range("S7").Value = 0.01361720793958475
$ws.Range("T7").Value = 0.01529925113856184
$ws.Range("I8").Value = 0.822180234441485
$ws.Range("J8").Value = 0.8526739017519405
$ws.Range("O8").Value = 0.2170932623988173
$ws.Range("P8").Value = 0.2351865659654651
$ws.Range("S8").Value = 0.1784897893747264
$ws.Range("T8").Value = 0.2005374468414133
$ws.Range("I9").Value = 0.822180234441485
$ws.Range("J9").Value = 0.8526739017519405
$ws.Range("M9").Value = 17.58286933333333
$ws.Range("N9").Value = 52.748608
$ws.Range("O9").Value = 0.2112656061941426
$ws.Range("P9").Value = 0.22887321273073
$ws.Range("Q9").Value = 208.6315287998578
$ws.Range("R9").Value = 1877.68375919872
$ws.Range("S9").Value = 0.1736984056301226
$ws.Range("T9").Value = 0.1951542153056135
$ws.Range("I10").Value = 0.822180234441485
$ws.Range("J10").Value = 0.8526739017519405
$ws.Range("M10").Value = 19.2082395
$ws.Range("N10").Value = 38.416479
$ws.Range("O10").Value = 0.2307951156866419
$ws.Range("P10").Value = 0.1666869194070983
$ws.Range("Q10").Value = 227.91754272106
$ws.Range("R10").Value = 1367.50525632636
$ws.Range("S10").Value = 0.1897551823231929
$ws.Range("T10").Value = 0.1421295859418618
$ws.Range("I11").Value = 0.822180234441485
$ws.Range("J11").Value = 0.8526739017519405
$ws.Range("M11").Value = 26.988955
$ws.Range("N11").Value = 80.966865
$ws.Range("O11").Value = 0.3242837008298742
$ws.Range("P11").Value = 0.3513106263825066
$ws.Range("Q11").Value = 320.2405043007333
$ws.Range("R11").Value = 2882.1645387066
$ws.Range("S11").Value = 0.2666196491738583
$ws.Range("T11").Value = 0.2995534025244901
$ws.Range("G12").Value = 0.37892
$ws.Range("H12").Value = 1.13676
$ws.Range("I12").Value = 0.02625574638939025
$ws.Range("J12").Value = 0.02722954016579943
$ws.Range("M12").Value = 1.378421333333333
$ws.Range("N12").Value = 4.135264
$ws.Range("O12").Value = 0.01656231489052403
$ws.Range("P12").Value = 0.01794267551419991
$ws.Range("Q12").Value = 0.5223114116266667
$ws.Range("R12").Value = 4.70080270464
$ws.Range("S12").Value = 0.0004348559393868207
$ws.Range("T12").Value = 0.0004885708035958125
$ws.Range("G13").Value = 0.37892
$ws.Range("H13").Value = 1.13676
$ws.Range("I13").Value = 0.02625574638939025
$ws.Range("J13").Value = 0.02722954016579943
$ws.Range("O13").Value = 0.2170932623988173
$ws.Range("P13").Value = 0.2351865659654651
$ws.Range("Q13").Value = 6.846282605280001
$ws.Range("R13").Value = 61.61654344752001
$ws.Range("S13").Value = 0.005699945640388697
$ws.Range("T13").Value = 0.00640402204441307
$ws.Range("G14").Value = 0.37892
$ws.Range("H14").Value = 1.13676
$ws.Range("I14").Value = 0.02625574638939025
$ws.Range("J14").Value = 0.02722954016579943
$ws.Range("M14").Value = 17.58286933333333
$ws.Range("N14").Value = 52.748608
$ws.Range("O14").Value = 0.2112656061941426
$ws.Range("P14").Value = 0.22887321273073
$ws.Range("Q14").Value = 6.662500847786665
$ws.Range("R14").Value = 59.96250763008
$ws.Range("S14").Value = 0.005546936177034201
$ws.Range("T14").Value = 0.00623211233892697
$ws.Range("G15").Value = 0.37892
$ws.Range("H15").Value = 1.13676
$ws.Range("I15").Value = 0.02625574638939025
$ws.Range("J15").Value = 0.02722954016579943
$ws.Range("M15").Value = 19.2082395
$ws.Range("N15").Value = 38.416479
$ws.Range("O15").Value = 0.2307951156866419
$ws.Range("P15").Value = 0.1666869194070983
$ws.Range("Q15").Value = 7.27838611134
$ws.Range("R15").Value = 43.67031666804
$ws.Range("S15").Value = 0.006059698025378453
$ws.Range("T15").Value = 0.004538808167108957
$ws.Range("G16").Value = 0.37892
$ws.Range("H16").Value = 1.13676
$ws.Range("I16").Value = 0.02625574638939025
$ws.Range("J16").Value = 0.02722954016579943
$ws.Range("M16").Value = 26.988955
$ws.Range("N16").Value = 80.966865
$ws.Range("O16").Value = 0.3242837008298742
$ws.Range("P16").Value = 0.3513106263825066
$ws.Range("Q16").Value = 10.2266548286
$ws.Range("R16").Value = 92.0398934574
$ws.Range("S16").Value = 0.008514310607202077
$ws.Range("T16").Value = 0.00956602681175462
$ws.Range("G17").Value = 1.548357
$ws.Range("H17").Value = 3.096714
$ws.Range("I17").Value = 0.1072872076222874
$ws.Range("J17").Value = 0.0741775733180209
$ws.Range("M17").Value = 1.378421333333333
$ws.Range("N17").Value = 4.135264
$ws.Range("O17").Value = 0.01656231489052403
$ws.Range("P17").Value = 0.01794267551419991
$ws.Range("Q17").Value = 2.134288320416001
$ws.Range("R17").Value = 12.805729922496
$ws.Range("S17").Value = 0.001776924516365353
$ws.Range("T17").Value = 0.001330944128476022
$ws.Range("G18").Value = 1.548357
$ws.Range("H18").Value = 3.096714
$ws.Range("I18").Value = 0.1072872076222874
$ws.Range("J18").Value = 0.0741775733180209
$ws.Range("O18").Value = 0.2170932623988173
$ws.Range("P18").Value = 0.2351865659654651
$ws.Range("Q18").Value = 27.97553466658801
$ws.Range("R18").Value = 167.8532079995281
$ws.Range("S18").Value = 0.02329132991638162
$ws.Range("T18").Value = 0.01744556874031685
$ws.Range("G19").Value = 1.548357
$ws.Range("H19").Value = 3.096714
$ws.Range("I19").Value = 0.1072872076222874
$ws.Range("J19").Value = 0.0741775733180209
$ws.Range("M19").Value = 17.58286933333333
$ws.Range("N19").Value = 52.748608
$ws.Range("O19").Value = 0.2112656061941426
$ws.Range("P19").Value = 0.22887321273073
$ws.Range("Q19").Value = 27.224558812352
$ws.Range("R19").Value = 163.347352874112
$ws.Range("S19").Value = 0.02266609695519937
$ws.Range("T19").Value = 0.01697725951786472
$ws.Range("G20").Value = 1.548357
$ws.Range("H20").Value = 3.096714
$ws.Range("I20").Value = 0.1072872076222874
$ws.Range("J20").Value = 0.0741775733180209
$ws.Range("M20").Value = 19.2082395
$ws.Range("N20").Value = 38.416479
$ws.Range("O20").Value = 0.2307951156866419
$ws.Range("P20").Value = 0.1666869194070983
$ws.Range("Q20").Value = 29.74121208750151
$ws.Range("R20").Value = 118.964848350006
$ws.Range("S20").Value = 0.02476136349488258
$ws.Range("T20").Value = 0.01236443118547508
$ws.Range("G21").Value = 1.548357
$ws.Range("H21").Value = 3.096714
$ws.Range("I21").Value = 0.1072872076222874
$ws.Range("J21").Value = 0.0741775733180209
$ws.Range("M21").Value = 26.988955
$ws.Range("N21").Value = 80.966865
$ws.Range("O21").Value = 0.3242837008298742
$ws.Range("P21").Value = 0.3513106263825066
$ws.Range("Q21").Value = 41.78853739693501
$ws.Range("R21").Value = 250.73122438161
$ws.Range("S21").Value = 0.03479149273945843
$ws.Range("T21").Value = 0.02605936974588823
$ws.Range("G22").Value = 0.16626
$ws.Range("H22").Value = 0.49878
$ws.Range("I22").Value = 0.01152032195371061
$ws.Range("J22").Value = 0.01194759671689489
$ws.Range("M22").Value = 1.378421333333333
$ws.Range("N22").Value = 4.135264
$ws.Range("O22").Value = 0.01656231489052403
$ws.Range("P22").Value = 0.01794267551419991
$ws.Range("Q22").Value = 0.22917633088
$ws.Range("R22").Value = 2.06258697792
$ws.Range("S22").Value = 0.0001908031998375721
$ws.Range("T22").Value = 0.0002143718510657653
$ws.Range("G23").Value = 0.16626
$ws.Range("H23").Value = 0.49878
$ws.Range("I23").Value = 0.01152032195371061
$ws.Range("J23").Value = 0.01194759671689489
$ws.Range("O23").Value = 0.2170932623988173
$ws.Range("P23").Value = 0.2351865659654651
$ws.Range("Q23").Value = 3.00396639384
$ws.Range("R23").Value = 27.03569754456
$ws.Range("S23").Value = 0.002500984276815752
$ws.Range("T23").Value = 0.002809914243386776
$ws.Range("G24").Value = 0.16626
$ws.Range("H24").Value = 0.49878
$ws.Range("I24").Value = 0.01152032195371061
$ws.Range("J24").Value = 0.01194759671689489
$ws.Range("M24").Value = 17.58286933333333
$ws.Range("N24").Value = 52.748608
$ws.Range("O24").Value = 0.2112656061941426
$ws.Range("P24").Value = 0.22887321273073
$ws.Range("Q24").Value = 2.923327855359999
$ws.Range("R24").Value = 26.30995069824
$ws.Range("S24").Value = 0.00243384780110236
$ws.Range("T24").Value = 0.002734484845006857
$ws.Range("G25").Value = 0.16626
$ws.Range("H25").Value = 0.49878
$ws.Range("I25").Value = 0.01152032195371061
$ws.Range("J25").Value = 0.01194759671689489
$ws.Range("M25").Value = 19.2082395
$ws.Range("N25").Value = 38.416479
$ws.Range("O25").Value = 0.2307951156866419
$ws.Range("P25").Value = 0.1666869194070983
$ws.Range("Q25").Value = 3.19356189927
$ws.Range("R25").Value = 19.16137139562
$ws.Range("S25").Value = 0.002658834038054
$ws.Range("T25").Value = 0.001991508091057572
$ws.Range("G26").Value = 0.16626
$ws.Range("H26").Value = 0.49878
$ws.Range("I26").Value = 0.01152032195371061
$ws.Range("J26").Value = 0.01194759671689489
$ws.Range("M26").Value = 26.988955
$ws.Range("N26").Value = 80.966865
$ws.Range("O26").Value = 0.3242837008298742
$ws.Range("P26").Value = 0.3513106263825066
$ws.Range("Q26").Value = 4.4871836583
$ws.Range("R26").Value = 40.3846529247
$ws.Range("S26").Value = 0.003735852637900922
$ws.Range("T26").Value = 0.004197317686377925
